# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list grew (older periods 1607-2104 were appended), and the
# table on Hoja1 (rows 16-74, E:G = Periodo Mora / Valor Mora / Salario Basico)
# is re-sorted from newest-period-first to oldest-period-first, with the
# Valor Mora / Salario Basico figures refreshed to match the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Ascending period list (oldest -> newest) that now fills rows 16-74.
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105"
)

# Refreshed "Valor Mora" (column F) for each of the same 59 rows, in row order.
$valorMora = @(
    24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,24640,
    24640,24640,24640,24640,24640,24640,24640,24640,
    31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,
    31249,31249,31249,31249,
    24640
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = 781242
}
